# Add team record (Wins / Losses / Ties) columns to the NYY 1994 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 38

# --- Header row (row 1): new columns AD, AE, AF -----------------------
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the look of the existing header cells (bold, thin box border,
# centered / top-aligned) used by every other column in row 1.
$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1         # xlContinuous
$headerRange.Borders.Weight = 2            # xlThin

# --- Data rows 2..38: constant team record for every player row -------
$ws.Range("AD2:AD$lastRow").Value = 70
$ws.Range("AE2:AE$lastRow").Value = 43
$ws.Range("AF2:AF$lastRow").Value = 0
